$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 3.25
$ws.Range("N2").Value = 9.75
$ws.Range("P2").Value = 15
$ws.Range("V2").Value = 16
$ws.Range("W2").Value = 90
$ws.Range("Y2").Value = 4.5
$ws.Range("Z2").Value = 5.6

# Row 9
$ws.Range("H9").Value = 3.6
$ws.Range("I9").Value = 6
$ws.Range("O9").Value = 6.5
$ws.Range("T9").Value = 7.5
$ws.Range("W9").Value = 81
$ws.Range("Z9").Value = 29
$ws.Range("AB9").Value = 67
$ws.Range("AE9").Value = 1.08
$ws.Range("AF9").Value = 7.5

# Row 14
$ws.Range("I14").Value = 2.75
$ws.Range("W14").Value = 67
$ws.Range("AI14").Value = 2.1
$ws.Range("AJ14").Value = 1.67

# Row 15
$ws.Range("G15").Value = 1.27
$ws.Range("H15").Value = 4.45
$ws.Range("I15").Value = 10.25
$ws.Range("J15").Value = 1.75
$ws.Range("K15").Value = 1.85
$ws.Range("N15").Value = 5
$ws.Range("O15").Value = 4.65
$ws.Range("P15").Value = 7.6
$ws.Range("Q15").Value = 6.1
$ws.Range("R15").Value = 10
$ws.Range("S15").Value = 28
$ws.Range("T15").Value = 10
$ws.Range("U15").Value = 8.25
$ws.Range("V15").Value = 21
$ws.Range("W15").Value = 110
$ws.Range("Y15").Value = 19
$ws.Range("Z15").Value = 60
$ws.Range("AA15").Value = 26
$ws.Range("AB15").Value = 250
$ws.Range("AC15").Value = 120
$ws.Range("AD15").Value = 90

# Row 17
$ws.Range("Y17").Value = 8
$ws.Range("AF17").Value = 8.5

# Row 18
$ws.Range("AC18").Value = 19
$ws.Range("AD18").Value = 26

# Row 19
$ws.Range("G19").Value = 1.36
$ws.Range("H19").Value = 3.8
$ws.Range("I19").Value = 8.5
$ws.Range("J19").Value = 2
$ws.Range("K19").Value = 1.8
$ws.Range("L19").Value = 1.4
$ws.Range("M19").Value = 2.75
$ws.Range("N19").Value = 6
$ws.Range("O19").Value = 6
$ws.Range("Q19").Value = 8.5
$ws.Range("R19").Value = 13
$ws.Range("S19").Value = 34
$ws.Range("T19").Value = 9
$ws.Range("V19").Value = 23
$ws.Range("W19").Value = 81
$ws.Range("Y19").Value = 19
$ws.Range("AD19").Value = 67
$ws.Range("AE19").Value = 1.07
$ws.Range("AF19").Value = 9
$ws.Range("AG19").Value = 1.3
$ws.Range("AH19").Value = 3.4
$ws.Range("AI19").Value = 2.25
$ws.Range("AJ19").Value = 1.57

# Row 20
$ws.Range("J20").Value = 2.08
$ws.Range("K20").Value = 1.73
$ws.Range("AF20").Value = 9

# Row 21
$ws.Range("G21").Value = 1.48
$ws.Range("H21").Value = 3.4
$ws.Range("I21").Value = 7
$ws.Range("J21").Value = 2.15
$ws.Range("K21").Value = 1.67
$ws.Range("L21").Value = 1.42
$ws.Range("M21").Value = 2.47
$ws.Range("N21").Value = 5.5
$ws.Range("O21").Value = 6
$ws.Range("P21").Value = 9
$ws.Range("Q21").Value = 10
$ws.Range("R21").Value = 15
$ws.Range("S21").Value = 34
$ws.Range("T21").Value = 7.5
$ws.Range("U21").Value = 7.5
$ws.Range("V21").Value = 21
$ws.Range("W21").Value = 81
$ws.Range("Y21").Value = 15
$ws.Range("Z21").Value = 34
$ws.Range("AA21").Value = 23
$ws.Range("AB21").Value = 81
$ws.Range("AC21").Value = 51
$ws.Range("AD21").Value = 67
$ws.Range("AG21").Value = 1.36
$ws.Range("AH21").Value = 2.67
$ws.Range("AI21").Value = 2.18
$ws.Range("AJ21").Value = 1.53

# Row 22
$ws.Range("X22").Value = 700

# Row 23
$ws.Range("G23").Value = 2.45
$ws.Range("I23").Value = 2.88
$ws.Range("O23").Value = 11
$ws.Range("P23").Value = 11
$ws.Range("Q23").Value = 26
$ws.Range("Y23").Value = 7

# Row 24
$ws.Range("G24").Value = 1.67
$ws.Range("H24").Value = 3.6
$ws.Range("I24").Value = 4.33
$ws.Range("J24").Value = 2.15
$ws.Range("K24").Value = 1.67
$ws.Range("L24").Value = 1.42
$ws.Range("M24").Value = 2.47
$ws.Range("N24").Value = 6
$ws.Range("O24").Value = 7.5
$ws.Range("P24").Value = 9
$ws.Range("Q24").Value = 13
$ws.Range("R24").Value = 15
$ws.Range("S24").Value = 34
$ws.Range("T24").Value = 9
$ws.Range("U24").Value = 7.5
$ws.Range("V24").Value = 21
$ws.Range("W24").Value = 67
$ws.Range("Y24").Value = 10
$ws.Range("Z24").Value = 23
$ws.Range("AA24").Value = 15
$ws.Range("AB24").Value = 51
$ws.Range("AC24").Value = 41
$ws.Range("AD24").Value = 41
$ws.Range("AG24").Value = 1.36
$ws.Range("AH24").Value = 2.67
$ws.Range("AI24").Value = 1.98
$ws.Range("AJ24").Value = 1.65

# Row 25
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 3
$ws.Range("I25").Value = 3.5
$ws.Range("J25").Value = 2.4
$ws.Range("K25").Value = 1.53
$ws.Range("L25").Value = 1.47
$ws.Range("M25").Value = 2.32
$ws.Range("N25").Value = 6
$ws.Range("O25").Value = 8.5
$ws.Range("P25").Value = 9.5
$ws.Range("Q25").Value = 19
$ws.Range("R25").Value = 21
$ws.Range("S25").Value = 34
$ws.Range("T25").Value = 7
$ws.Range("U25").Value = 6
$ws.Range("V25").Value = 19
$ws.Range("W25").Value = 67
$ws.Range("Y25").Value = 8.5
$ws.Range("Z25").Value = 17
$ws.Range("AA25").Value = 13
$ws.Range("AB25").Value = 41
$ws.Range("AC25").Value = 34
$ws.Range("AD25").Value = 41
$ws.Range("AG25").Value = 1.42
$ws.Range("AH25").Value = 2.47
$ws.Range("AI25").Value = 1.93
$ws.Range("AJ25").Value = 1.7
